$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 18740
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 18740
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 18740
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -19326
# Row 108
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# Row 138
$ws.Range("H138").Value = 5596.5713
$ws.Range("I138").Value = 5315.8
$ws.Range("J138").Value = 6298.5
$ws.Range("K138").Value = 15947.4
$ws.Range("L138").Value = 18895.5
$ws.Range("M138").Value = -10807.4
$ws.Range("N138").Value = -29175.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4016.1428
$ws.Range("I45").Value = 3576.3333
$ws.Range("J45").Value = 4807.8
$ws.Range("K45").Value = 3576.3333
$ws.Range("L45").Value = 4807.8
$ws.Range("M45").Value = -3199.3333
$ws.Range("N45").Value = -5561.8
# Row 80
$ws.Range("H80").Value = 79999
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 79999
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2387.6875
$ws.Range("I105").Value = 2266.0833
$ws.Range("J105").Value = 2752.5
$ws.Range("K105").Value = 2266.0833
$ws.Range("L105").Value = 2752.5
$ws.Range("M105").Value = -519.0832999999998
$ws.Range("N105").Value = -6246.5
# Row 107
$ws.Range("H107").Value = 2094.348
$ws.Range("I107").Value = 1272.625
$ws.Range("J107").Value = 3972.5715
$ws.Range("K107").Value = 1272.625
$ws.Range("L107").Value = 3972.5715
$ws.Range("M107").Value = 647.375
$ws.Range("N107").Value = -7812.5715

$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 1066.7778
$ws.Range("I6").Value = 667
$ws.Range("J6").Value = 1266.6666
$ws.Range("K6").Value = 667
$ws.Range("L6").Value = 1266.6666
$ws.Range("M6").Value = -554
$ws.Range("N6").Value = -1492.6666
# Row 99
$ws.Range("H99").Value = 3066.0667
$ws.Range("I99").Value = 2779.1
$ws.Range("J99").Value = 3640
$ws.Range("K99").Value = 2779.1
$ws.Range("L99").Value = 3640
$ws.Range("M99").Value = -1281.1
$ws.Range("N99").Value = -6636
# Row 126
$ws.Range("H126").Value = 3066.0667
$ws.Range("I126").Value = 2779.1
$ws.Range("J126").Value = 3640
$ws.Range("K126").Value = 8337.299999999999
$ws.Range("L126").Value = 10920
$ws.Range("M126").Value = -5867.299999999999
$ws.Range("N126").Value = -15860

$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 1140.7
$ws.Range("I46").Value = 1334.6666
$ws.Range("J46").Value = 849.75
$ws.Range("K46").Value = 4003.9998
$ws.Range("L46").Value = 2549.25
$ws.Range("M46").Value = -3912.9998
$ws.Range("N46").Value = -2731.25
# Row 55
$ws.Range("H55").Value = 3810.818
$ws.Range("I55").Value = 560.2857
$ws.Range("J55").Value = 9499.25
$ws.Range("K55").Value = 1680.8571
$ws.Range("L55").Value = 28497.75
$ws.Range("M55").Value = -1503.8571
$ws.Range("N55").Value = -28851.75
# Row 56
$ws.Range("H56").Value = 5236
$ws.Range("I56").Value = 5236
$ws.Range("K56").Value = 5236
$ws.Range("M56").Value = -4706
# Row 104
$ws.Range("H104").Value = 1499
$ws.Range("I104").Value = 2499
$ws.Range("J104").Value = 999
$ws.Range("K104").Value = 7497
$ws.Range("L104").Value = 2997
$ws.Range("M104").Value = -4876
$ws.Range("N104").Value = -8239
# Row 121
$ws.Range("H121").Value = 599
$ws.Range("I121").Value = 491.2
$ws.Range("J121").Value = 778.6667
$ws.Range("K121").Value = 1473.6
$ws.Range("L121").Value = 2336.0001
$ws.Range("M121").Value = -163.5999999999999
$ws.Range("N121").Value = -4956.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 606.5
$ws.Range("I97").Value = 638.08
$ws.Range("J97").Value = 493.7143
$ws.Range("K97").Value = 638.08
$ws.Range("L97").Value = 493.7143
$ws.Range("M97").Value = -142.08
$ws.Range("N97").Value = -1485.7143

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 13431.833
$ws.Range("I46").Value = 31374.5
$ws.Range("J46").Value = 4460.5
$ws.Range("K46").Value = 31374.5
$ws.Range("L46").Value = 4460.5
$ws.Range("M46").Value = -31186.5
$ws.Range("N46").Value = -4836.5
# Row 55
$ws.Range("H55").Value = 938.4167
$ws.Range("I55").Value = 872.4286
$ws.Range("K55").Value = 872.4286
$ws.Range("M55").Value = -699.4286
# Row 82
$ws.Range("H82").Value = 2053.7083
$ws.Range("I82").Value = 1776.8462
$ws.Range("J82").Value = 2380.9092
$ws.Range("K82").Value = 1776.8462
$ws.Range("L82").Value = 2380.9092
$ws.Range("M82").Value = -1415.8462
$ws.Range("N82").Value = -3102.9092
# Row 85
$ws.Range("H85").Value = 2053.7083
$ws.Range("I85").Value = 1776.8462
$ws.Range("J85").Value = 2380.9092
$ws.Range("K85").Value = 1776.8462
$ws.Range("L85").Value = 2380.9092
$ws.Range("M85").Value = -528.8462
$ws.Range("N85").Value = -4876.9092

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 503332.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 503332.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 503332.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -504580.5
# Row 65
$ws.Range("H65").Value = 503332.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 503332.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 2516662.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -2522902.5
# Row 107
$ws.Range("H107").Value = 1456.5714
$ws.Range("I107").Value = 713.7619
$ws.Range("J107").Value = 3685
$ws.Range("K107").Value = 2141.2857
$ws.Range("L107").Value = 11055
$ws.Range("M107").Value = -221.2856999999999
$ws.Range("N107").Value = -14895
# Row 113
$ws.Range("H113").Value = 1299.421
$ws.Range("I113").Value = 513.7143
$ws.Range("J113").Value = 3499.4
$ws.Range("K113").Value = 1541.1429
$ws.Range("L113").Value = 10498.2
$ws.Range("M113").Value = 628.8571000000002
$ws.Range("N113").Value = -14838.2
